# Assignment : Checking visibility of pagelinks on Account page
#
# 1. Rename "Sheet2" -> "pagelinks" and populate it with the account page's
#    side-navigation links (header "links" + 13 menu items).
# 2. Make "pagelinks" the active tab / selected sheet, with C17 selected.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("register")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Name = "pagelinks"

# Menu items, in row order starting at A2 (A1 holds the "links" header and is
# written afterwards so its shared-string entry lands last, matching the
# order new strings were appended to the workbook).
$menuItems = @(
    "My Account",
    "Edit Account",
    "Password",
    "Address Book",
    "Wish List",
    "Order History",
    "Downloads",
    "Recurring payments",
    "Reward Points",
    "Returns",
    "Transactions",
    "Newsletter",
    "Logout"
)

for ($i = 0; $i -lt $menuItems.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $menuItems[$i]
}

$ws2.Range("A1").Value = "links"

# Give the header cell the same highlighted fill used by the "register"
# sheet's header row.
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Column A is sized to fit its widest entry ("Recurring payments").
$ws2.Columns.Item(1).ColumnWidth = 18.85546875

# Switch focus to the new sheet and select C17, matching the saved view.
$ws2.Activate() | Out-Null
$ws2.Range("C17").Select() | Out-Null
